$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet gets a new column inserted before column N,
# shifting the old N/O/P (Late / heading / Outstanding) columns one place
# to the right, and becomes the active/selected sheet+cell.
$ws = $wb.Worksheets.Item("Repayment schedule")

$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab (this also clears the
# previously active "Summary" sheet's tabSelected flag) and move the
# selection to the new active cell.
$ws.Activate()
$ws.Range("R9").Select() | Out-Null
